$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ProductBacklog / SprintBacklog status updates ---------------------------------
# Row 7  (ID 5):  User Story cell re-typed (font refresh) + Story Status IP -> D
$ws.Range("B7").Font.Name = "Arial"
$ws.Range("B7").Font.Size = 12
$ws.Range("B7").Font.ThemeColor = 1

$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "D"

# Row 9  (ID 7):  Story Status W -> D
$ws.Range("E2").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "D"

# Row 11 (ID 9):  Story Status W -> D
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "D"

# Row 12 (ID 10): User Story cell re-typed (font refresh) + Story Status W -> IP
$ws.Range("B12").Font.Name = "Arial"
$ws.Range("B12").Font.Size = 12
$ws.Range("B12").Font.ThemeColor = 1

$ws.Range("E8").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = "IP"

# Row 14 (ID 12): Story Status W -> IP
$ws.Range("E8").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = "IP"

# --- View state: zoom + active selection -------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("E14").Select()
